$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column E ("MP") holds ratio values for rows 2-283 (row 1 is the header).
# The edit replaces every value with its reciprocal (1/x).
$lastRow = 283
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 5)
    $v = $cell.Value2
    $cell.Value2 = 1 / $v
}
